# Add correlated PARAM5 ~ TruncatedNormal(3, 1, 1, 5) and PARAM6 ~ Uniform(0, 1)
# to design_input_background.xlsx, with a new "corr0" correlation sheet instead
# of sampling from multivariate_normal directly in tests.

$wb = $excel.ActiveWorkbook

# --- 1. Insert a new correlation sheet "corr0" right before "corr1" --------
$corr1 = $wb.Worksheets.Item("corr1")
$corr0 = $wb.Worksheets.Add($corr1)
$corr0.Name = "corr0"

$corr0.Range("B1").Value = "PARAM5"
$corr0.Range("C1").Value = "PARAM6"
$corr0.Range("A2").Value = "PARAM5"
$corr0.Range("B2").Value = 1
$corr0.Range("A3").Value = "PARAM6"
$corr0.Range("B3").Value = 0.8
$corr0.Range("C3").Value = 1

# --- 2. designinput: sens6 / PARAM5 & PARAM6 now use the corr0 sheet -------
$designinput = $wb.Worksheets.Item("designinput")
$designinput.Range("B9").Value = 500
$designinput.Range("O9").Value = "corr0"
$designinput.Range("O10").Value = "corr0"

# Selection bookkeeping to mirror the saved workbook state.
$designinput.Range("B10").Select() | Out-Null
$corr0.Range("C8").Select() | Out-Null
$corr0.Activate() | Out-Null
